# Atualização automática de PEDRO_OSORIO.xlsx
#
# 1) Rename "Paineis DARQ"            -> "PAINEIS DARQ"
# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3) Delete the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete" prompt when removing a sheet
$excel.DisplayAlerts = $false | Out-Null

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

$excel.DisplayAlerts = $true | Out-Null
